$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column J: "Rule For" ------------------------------------------------
# Header cell J1 — same bold/no-fill "header" look as A1:I1, but without the
# bottom border that the other header cells have.
$ws.Range("J1").Value = "Rule For"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("J1").Borders.LineStyle = -4142        # xlLineStyleNone (drop the header border)

# Data cells J2:J20 — same look as the rest of the data rows (e.g. column I).
$ws.Range("J2").Value = "Accounting"
$ws.Range("J3").Value = "Reporting"
for ($r = 4; $r -le 20; $r++) {
    $ws.Cells.Item($r, 10).Value = "Accounting"
}

$ws.Range("I2").Copy() | Out-Null
$ws.Range("J2:J20").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Selection ----------------------------------------------------------
$ws.Range("J4:J20").Select() | Out-Null
